# Edit: "fixed bounds not updating properly in S.bo"
# Applies the value/format/selection changes recorded in the target diff
# across the four sheets of paramTables.xlsx: Vmax, Ks,
# "Initial & Flow Concentrations", and "Other Parameters".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Vmax"
# ---------------------------------------------------------------------
$vmax = $wb.Worksheets.Item("Vmax")

# Row 7 (Glucose bounds) - S.bo lower bound corrected (was 10x too high),
# B/C (B.th / E.re) reset to a flat 0.1.
$vmax.Range("A7").Value = 0.89922067541464101
$vmax.Range("B7").Value = 0.1
$vmax.Range("C7").Value = 0.1

# Row 10 (Phosphate bounds): B/C/D used to hold a stale shared value;
# move that stale value out to the new I/J/K columns and make B/C/D
# mirror column A like the rest of the row.
$vmax.Range("I10").Value = 0.12635382862630201
$vmax.Range("J10").Value = 0.12635382862630201
$vmax.Range("K10").Value = 0.12635382862630201
$vmax.Range("B10").Value = 0.21999999999999997
$vmax.Range("C10").Value = 0.21999999999999997
$vmax.Range("D10").Value = 0.21999999999999997
# C10 loses its Arial override (now matches the plain default style of B10/D10).
$vmax.Range("B10").Copy()
$vmax.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
# J10 picks up the Arial override that C10 used to have.
$vmax.Range("C9").Copy()
$vmax.Range("J10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$vmax.Range("J10").Value = 0.12635382862630201

# B14 gains the Arial style used elsewhere in the column (value unchanged).
$vmax.Range("A14").Copy()
$vmax.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 19-20 (Tryptophan / Methionine): B/C now mirror column A instead of 0.
$vmax.Range("B19").Value = 0.23599999999999999
$vmax.Range("C19").Value = 0.23599999999999999
$vmax.Range("B20").Value = 0.156
$vmax.Range("C20").Value = 0.156

$vmax.Activate()
$vmax.Range("B20:C20").Select()

# ---------------------------------------------------------------------
# Sheet "Ks"
# ---------------------------------------------------------------------
$ks = $wb.Worksheets.Item("Ks")

# Rows 19-20: B/C now mirror column A instead of 0.
$ks.Range("B19").Value = 0.00089999999999999998
$ks.Range("C19").Value = 0.00089999999999999998
$ks.Range("B20").Value = 0.0022699999999999999
$ks.Range("C20").Value = 0.0022699999999999999

$ks.Activate()
$ks.Range("B57").Select()

# ---------------------------------------------------------------------
# Sheet "Initial & Flow Concentrations"
# ---------------------------------------------------------------------
$ifc = $wb.Worksheets.Item("Initial & Flow Concentrations")

# S.bo flow concentration bound fix (was 50, should be 5) - this is the
# headline fix referenced by the commit message.
$ifc.Range("B1").Value = 5

# Cancer / Colon biomass initial conditions cleared; the old Colon
# biomass initial value (1.7) is preserved off to the side in column F.
$ifc.Range("A5").Value = 0
$ifc.Range("F6").Value = 1.7
$ifc.Range("A6").Value = 0

# Glucose / Phosphate / Ammonium initial conditions and flow
# concentrations retuned.
$ifc.Range("A7").Value = 10
$ifc.Range("B8").Value = 500
$ifc.Range("B9").Value = 50
$ifc.Range("A10").Value = 10
$ifc.Range("B10").Value = 200
$ifc.Range("A11").Value = 20
$ifc.Range("B11").Value = 100

# E7 (stray empty Arial-styled placeholder) removed entirely.
$ifc.Range("E7").Clear()

# Tryptophan flow concentration: replace the stale "=D19*10" formula
# (D19 is empty, so it always evaluated to 0) with a literal value that
# mirrors the initial condition, like every other row.
$ifc.Range("B19").Value = 0.47199999999999998

# MFalpha2 flow concentration reset to 0.
$ifc.Range("B29").Value = 0

$ifc.Activate()
$ifc.Range("K20").Select()

# ---------------------------------------------------------------------
# Sheet "Other Parameters"
# ---------------------------------------------------------------------
$other = $wb.Worksheets.Item("Other Parameters")

# MFalpha2 anti-cancer-secretion threshold raised from 0.001 to 0.01.
$other.Range("A4").Value = 0.01

$other.Activate()
$other.Range("A5").Select()

# ---------------------------------------------------------------------
# Leave "Initial & Flow Concentrations" as the active tab/selection,
# matching the saved workbook view.
# ---------------------------------------------------------------------
$ifc.Activate()
$ifc.Range("K20").Select()
